# "added themeConfig for theme, styleConfig for each component"
#
# Both heading components move from slide-specific, hard-coded local
# formatting to the shared theme/style configuration:
#   - the box collapses to the full-width banner placement (0,0)-(540pt,0pt)
#     driven by themeConfig instead of its own hand-placed xfrm
#   - the box fill becomes the theme's white background instead of a
#     per-shape colour
#   - the run text recolors to the brand blue (305597) and switches from
#     Arial to the new Roboto brand typeface
#   - the local size/bold/alignment paragraph & run overrides are dropped
#     so the component falls back to the shared styleConfig defaults
#     (left aligned, 18pt, regular weight)

$p = $ppt.ActivePresentation

# EMU -> points conversion used throughout: 1 pt = 12700 EMU
# 457200/228600/etc EMU offsets collapse to 0; 6858000 EMU (540 pt) is the
# full 10in slide width; cy collapses to 0.

# ---------------------------------------------------------------------
# Slide 1 : "Welcome to Our Presentation!"
# ---------------------------------------------------------------------
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)

# p:spPr/a:xfrm : off (457200,228600) ext (5486400,1371600) -> off (0,0) ext (6858000,0)
$sh1.Left   = 0
$sh1.Top    = 0
$sh1.Width  = 540
$sh1.Height = 0

# p:spPr/a:solidFill : 305597 -> FFFFFF
$sh1.Fill.ForeColor.RGB = 16777215

$tr1 = $sh1.TextFrame.TextRange

# a:pPr : algn="ctr" override removed -> falls back to the default (left)
$tr1.ParagraphFormat.Alignment = 1

# a:rPr : sz="4800" b="1" overrides removed -> falls back to the default
# (18pt, regular weight); fill recolors FFFFFF -> 305597
$tr1.Font.Size = 18
$tr1.Font.Bold = $false
$tr1.Font.Color.RGB = 9917744

# a:latin / a:ea / a:cs : Arial -> Roboto
$tr1.Font.Name             = "Roboto"
$tr1.Font.NameFarEast      = "Roboto"
$tr1.Font.NameComplexScript = "Roboto"

# ---------------------------------------------------------------------
# Slide 2 : "First Heading"
# ---------------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)

# p:spPr/a:xfrm : off (457200,1600200) ext (5486400,457200) -> off (0,0) ext (6858000,0)
$sh2.Left   = 0
$sh2.Top    = 0
$sh2.Width  = 540
$sh2.Height = 0

# p:spPr/a:solidFill : 000000 -> FFFFFF
$sh2.Fill.ForeColor.RGB = 16777215

$tr2 = $sh2.TextFrame.TextRange

# a:pPr : algn="l" override removed -- already the default alignment, so
# no property change is needed here.

# a:rPr : sz="3200" b="1" overrides removed -> falls back to the default
# (18pt, regular weight); fill color (305597) is unchanged
$tr2.Font.Size = 18
$tr2.Font.Bold = $false

# a:latin / a:ea / a:cs : Arial -> Roboto
$tr2.Font.Name              = "Roboto"
$tr2.Font.NameFarEast       = "Roboto"
$tr2.Font.NameComplexScript = "Roboto"
